# Add a "Skill Description" column (full skill name) after SkillCode.
# For generic SFIA attribute groups (Autonomy, Influence, Complexity, Knowledge)
# the description repeats the group name; for specific skill codes it is the
# human friendly full name of the skill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map SkillCode (column A) -> Skill Description (new column B)
$friendly = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "RESC"       = "Resourcing"
    "LEDA"       = "Competency assessment"
    "SUPP"       = "Supplier management"
    "ITCM"       = "Contract management"
}

# Insert a new blank column before column B; existing B,C,D shift to C,D,E
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Skill Description"

# Determine the last used row on the sheet (data originally ran through row 25)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    if ($code -and $friendly.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $friendly[$code]
    }
}
